# Generate Report for Archive
# Two files (05f1f749-082e-4df3-a5d8-b4357d823046.md and
# 86d63155-4055-410d-b618-c1071449da11.md) have moved from
# "Ready for handoff" to "In Translation" status. Update the
# Overview summary sheet as well as the per-locale (zh-cn, de-de)
# detail sheets to reflect the new status.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E3").Value = "In Translation"
$ws.Range("F3").Value = "In Translation"
$ws.Range("E4").Value = "In Translation"
$ws.Range("F4").Value = "In Translation"

$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C3").Value = "In Translation"
$ws.Range("C4").Value = "In Translation"

$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C3").Value = "In Translation"
$ws.Range("C4").Value = "In Translation"
